# Update column F (dSF) values for specific rows, per repull of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -1
$ws.Range("F3").Value = -2
$ws.Range("F5").Value = -2
$ws.Range("F7").Value = -2
$ws.Range("F15").Value = -5
